$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 840, shifting rows 840:881 down to 841:882.
$ws.Rows.Item(840).Insert()

# Populate the newly inserted row 840 with the new entry.
# Column A holds the date as plain text elsewhere in the sheet (t="inlineStr"),
# so briefly force a text number format before assigning the value -
# otherwise Excel auto-converts the "yyyy/mm/dd"-looking string into a date
# serial number. Restore the cell to the default "Normal" style afterwards
# so no stray number format is left applied to the cell.
$ws.Cells.Item(840, 1).NumberFormat = "@"
$ws.Cells.Item(840, 1).Value = "2026/02/23"
$ws.Cells.Item(840, 1).Style = "Normal"

$ws.Cells.Item(840, 2).Value = "月"
$ws.Cells.Item(840, 3).Value = 13
$ws.Cells.Item(840, 4).Value = 28
